# Update gh-pages to output generated at 456a3b4
# Applies updated "想去人数" (F column) counts across the 展览, 演出 and 全部类型 sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 53
$ws1.Range("F4").Value = 1098
$ws1.Range("F5").Value = 368
$ws1.Range("F6").Value = 614
$ws1.Range("F7").Value = 588
$ws1.Range("F8").Value = 1511
$ws1.Range("F10").Value = 1414
$ws1.Range("F11").Value = 3051
$ws1.Range("F12").Value = 555
$ws1.Range("F13").Value = 1720
$ws1.Range("F14").Value = 1781
$ws1.Range("F17").Value = 1442
$ws1.Range("F19").Value = 71
$ws1.Range("F21").Value = 386
$ws1.Range("F23").Value = 56
$ws1.Range("F24").Value = 4637
$ws1.Range("F28").Value = 27
$ws1.Range("F29").Value = 78

# --- Sheet "演出" (performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F6").Value = 59
$ws2.Range("F9").Value = 46
$ws2.Range("F12").Value = 29
$ws2.Range("F13").Value = 96

# --- Sheet "全部类型" (all types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 53
$ws4.Range("F9").Value = 59
$ws4.Range("F13").Value = 46
$ws4.Range("F15").Value = 1098
$ws4.Range("F16").Value = 368
$ws4.Range("F17").Value = 614
$ws4.Range("F18").Value = 588
$ws4.Range("F19").Value = 1511
$ws4.Range("F21").Value = 1414
$ws4.Range("F22").Value = 3051
$ws4.Range("F23").Value = 555
$ws4.Range("F24").Value = 1720
$ws4.Range("F25").Value = 1781
$ws4.Range("F28").Value = 1442
$ws4.Range("F30").Value = 71
$ws4.Range("F32").Value = 29
$ws4.Range("F34").Value = 386
$ws4.Range("F36").Value = 56
$ws4.Range("F37").Value = 4637
$ws4.Range("F41").Value = 96
$ws4.Range("F43").Value = 27
$ws4.Range("F44").Value = 78
